# Update "想去人数" (interest count) figures in column F for the
# exhibition ("展览") and "全部类型" worksheets, reflecting a refreshed
# scrape of the source data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Worksheets index 1 / sheet1.xml) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value  = 533
$ws1.Range("F7").Value  = 1617
$ws1.Range("F9").Value  = 16
$ws1.Range("F10").Value = 1409
$ws1.Range("F12").Value = 21
$ws1.Range("F13").Value = 358
$ws1.Range("F14").Value = 241
$ws1.Range("F15").Value = 182
$ws1.Range("F16").Value = 5
$ws1.Range("F18").Value = 14
$ws1.Range("F19").Value = 254
$ws1.Range("F20").Value = 138
$ws1.Range("F21").Value = 208
$ws1.Range("F22").Value = 192

# --- Sheet 4: 全部类型 (Worksheets index 4 / sheet4.xml) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F6").Value  = 533
$ws4.Range("F7").Value  = 1617
$ws4.Range("F10").Value = 16
$ws4.Range("F11").Value = 1409
$ws4.Range("F13").Value = 21
$ws4.Range("F14").Value = 358
$ws4.Range("F15").Value = 241
$ws4.Range("F16").Value = 182
$ws4.Range("F17").Value = 5
$ws4.Range("F19").Value = 14
$ws4.Range("F20").Value = 254
$ws4.Range("F21").Value = 138
$ws4.Range("F22").Value = 208
$ws4.Range("F23").Value = 192
